# Actualización automática 2025-09-18 09:12:30
#
# Updates a handful of monthly-sales figures (and the downstream totals /
# "X de 30" counters / compliance percentages that are derived from them)
# for three advisor-client rows in "LOZANO MOLINA TITO" across the three
# worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" — per-category sales by client
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# PIEDRA SINTERIZADA (col L) for three clients
$wsGrupo.Range("L6").Value  = 655.34                # ARELLANO CEDEÑO DANNY MARCELO
$wsGrupo.Range("L19").Value = 556.8099999999999      # MATERIALES PARA DECORACION DECORCASA CIA. LTDA.
$wsGrupo.Range("L22").Value = 565.25                # RENOVA&DISEÑA S.A.

# 240X80 PORCELANATO (col D) for ROCA REYNA PAUL DAVID
$wsGrupo.Range("D24").Value = 1900.8

# "X de 30" summary row (row 32): recompute the two counters that moved
# because D24 and L22 went from 0 to a positive value.
$wsGrupo.Range("D32").Value = "1 de 30"
$wsGrupo.Range("L32").Value = "5 de 30"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" — septiembre column (F) by client, plus totals
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F6").Value  = 655.34               # ARELLANO CEDEÑO DANNY MARCELO
$wsMensual.Range("F19").Value = 556.8099999999999     # MATERIALES PARA DECORACION DECORCASA CIA. LTDA.
$wsMensual.Range("F22").Value = 565.25               # RENOVA&DISEÑA S.A.
$wsMensual.Range("F24").Value = 2868.69              # ROCA REYNA PAUL DAVID

# Column total for septiembre (row 32)
$wsMensual.Range("F32").Value = 8758.23

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" — VENTA / POR CUMPLIR / CUMPLIMIENTO
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 1900.8
$wsCumpl.Range("E3").Value = 681.61380675037
$wsCumpl.Range("F3").Value = 0.7360555442475379

# Row 11: PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 4225.13
$wsCumpl.Range("E11").Value = 1619.31916370549
$wsCumpl.Range("F11").Value = 0.7229304048426677

# Row 15: TOTAL
$wsCumpl.Range("D15").Value = 8832.52
$wsCumpl.Range("E15").Value = 22875.23990313501
$wsCumpl.Range("F15").Value = 0.2785602018869429
